# Fix the out-of-order "Linha BA"/"Linha BB" placeholder values in column B
# so the rows read Linha B2 / Linha B3 / Linha B4 in order (matches column A's
# Linha A2 / Linha A3 / Linha A4 pattern).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Linha B2"
$ws.Range("B3").Value = "Linha B3"

# Leave the selection on the last cell touched while reviewing the fix.
$ws.Range("B3").Select() | Out-Null

# Zoom in to double-check the (previously empty-looking / mismatched) cells.
$excel.ActiveWindow.Zoom = 235
